# "boot camp" -> "bootcamp" on the Conclusion feedback slide.
#
# The original run's text "What did you think of this boot camp?" gets
# split into three runs so that "bootcamp" (the corrected spelling) is
# isolated, matching how PowerPoint splits a run when you retype/replace
# a word in the middle of a sentence.

$p = $ppt.ActivePresentation

$target = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text.IndexOf("boot camp") -ge 0) {
                $target = $shape.TextFrame.TextRange
                $targetShape = $shape
            }
        }
    }
}

$fullText = $target.Text
$searchText = "boot camp"
$idx = $fullText.IndexOf($searchText)

# 1-based character position of the first char of "boot camp"
$startPos = $idx + 1

# Replacing just the "boot camp" span (not the whole run/paragraph) makes
# PowerPoint split the surrounding run into "...this " + "bootcamp" + "?"
# runs, leaving the rest of the sentence/formatting untouched.
$wholeMatch = $target.Characters($startPos, $searchText.Length)
$wholeMatch.Text = "bootcamp"

